$wb = $excel.ActiveWorkbook

# Select Tourism sheet and fill row 3 with the new effort entry
$ws = $wb.Worksheets.Item("Tourism")
$ws.Range("A3").Value = "27.10.2016"
$ws.Range("B3").Value = "Arpan Kar"
$ws.Range("C3").Value = "Operation Management"
$ws.Range("D3").Value = 0.5
$ws.Range("J3").Value = "Calculating sales commision, and how it will fit"

# Update the active selection on Tourism to A4
$ws.Range("A4").Select()

# Recalculate so dependent formulas across sheets update
$excel.Calculate()

# Make "Cover" the active sheet / selected tab, and unselect School
$wb.Worksheets.Item("Cover").Select()
